# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to match the newly scraped totals.

$wb = $excel.ActiveWorkbook

$changes = @{
    2  = 7238
    5  = 22
    6  = 566
    7  = 187
    12 = 220
    14 = 462
    16 = 1860
    18 = 41
    19 = 3774
    21 = 252
    26 = 2445
    28 = 303
    30 = 6
    31 = 42
    33 = 18
    38 = 1462
    39 = 156
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $changes.Keys) {
    $ws1.Range("F$row").Value = $changes[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $changes.Keys) {
    # 全部类型 carries one extra row (the 演出 sheet's single event, inserted
    # at row 6), so every 展览 row from 6 onward is shifted down by one.
    if ($row -ge 6) {
        $targetRow = $row + 1
    } else {
        $targetRow = $row
    }
    $ws4.Range("F$targetRow").Value = $changes[$row]
}
